# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" on all sheets
# - Narrow the (now shorter) status column(s) to reflect the new text width

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the "Ready for handoff" status values to "In Translation"
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Re-size the Status columns now that the text is shorter
$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511

$zhcn.Columns.Item(3).ColumnWidth = 13.4101845877511
$dede.Columns.Item(3).ColumnWidth = 13.4101845877511
